$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Update row 46: hours 4 -> 5.5, and text for D46 (combine the two log entries)
$ws.Range("B46").Value = 5.5
$ws.Range("D46").Value = "Indie Project: Eliminated the NPE!!!`r`nWeek 7/Indie Project: drafted servlet and jsp for admin to use, to demo authentication.  Not working yet, probably due to path issues."
$ws.Rows(46).RowHeight = 45

# Delete rows 48:50 (the blank placeholder rows + the old "Friday 9pm" entry),
# shifting everything below up by 3 rows.
$ws.Rows("48:50").Delete()

# Update the selection to match the post-edit state (rows 48:50 selected prior to delete)
$ws.Range("A48:XFD50").Select()
